# Apply the cryptos-list price/volume refresh (GitHub Actions data pull).
# Text-ish "Price" cells that look numeric (e.g. "1.00", "0.130") are written
# with a leading apostrophe so Excel keeps them as literal text instead of
# collapsing them to a Double and losing trailing zeros / exact formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.511.29"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "3.714.91"
$ws.Range("E3").Value = "  -5.59%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'597.27"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "'167.37"
$ws.Range("E6").Value = "  -4.57%  "
$ws.Range("D7").Value = "3.715.17"
$ws.Range("E7").Value = "  -5.56%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  -4.12%  "
$ws.Range("D12").Value = "'0.464"
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("D13").Value = "'38.00"
$ws.Range("E13").Value = "  -5.41%  "
$ws.Range("D14").Value = "'0.0000243"
$ws.Range("E14").Value = "  -5.17%  "
$ws.Range("D15").Value = "4.331.73"
$ws.Range("E15").Value = "  -5.25%  "
$ws.Range("D16").Value = "3.718.02"
$ws.Range("E16").Value = "  -5.30%  "
$ws.Range("D17").Value = "67.548.85"
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.25"
$ws.Range("E18").Value = "  -3.65%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.56"
$ws.Range("E19").Value = "  +4.67%  "
$ws.Range("D20").Value = "'0.115"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").Value = "'487.09"
$ws.Range("E21").Value = "  -4.33%  "
$ws.Range("D22").Value = "'9.31"
$ws.Range("E22").Value = "  -4.40%  "
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D24").Value = "'85.52"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("E25").Value = "  -6.45%  "
$ws.Range("D26").Value = "'0.0000137"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "'12.22"
$ws.Range("E27").Value = "  -3.82%  "
$ws.Range("D28").Value = "'10.13"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'2.94"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("E31").Value = "  -8.99%  "
$ws.Range("D32").Value = "'7.73"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "'31.55"
$ws.Range("E33").Value = "  -7.57%  "
$ws.Range("D34").Value = "3.854.80"
$ws.Range("E34").Value = "  -5.58%  "
$ws.Range("E35").Value = "  -4.81%  "
$ws.Range("D36").Value = "3.658.56"
$ws.Range("E36").Value = "  -5.25%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("D39").Value = "'5.82"
$ws.Range("E39").Value = "  -6.31%  "
$ws.Range("D40").Value = "'0.130"
$ws.Range("E40").Value = "  -7.72%  "
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("D42").Value = "'429.65"
$ws.Range("E42").Value = "  -8.57%  "
$ws.Range("D43").Value = "'48.69"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("E44").Value = "  -5.51%  "
$ws.Range("D45").Value = "'2.81"
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("D46").Value = "'8.44"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -5.29%  "
$ws.Range("D49").Value = "'141.58"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").Value = "2.753.75"
$ws.Range("E50").Value = "  -7.18%  "
$ws.Range("D51").Value = "'0.0349"
$ws.Range("E51").Value = "  -4.36%  "
